# Updated cryptos list on Mon Feb 26 11:44:48 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.098.53'
$ws.Range('E2').Value = '  -1.06%  '
$ws.Range('D3').Value = '3.054.62'
$ws.Range('E3').Value = '  +0.74%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = "'391.05"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.17%  '
$ws.Range('D6').Value = "'100.85"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.80%  '
$ws.Range('E7').Value = '  -2.23%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').Value = "'0.578"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.94%  '
$ws.Range('D10').Value = "'36.59"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.32%  '
$ws.Range('E11').Value = '  +0.37%  '
$ws.Range('D12').Value = "'0.0845"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.89%  '
$ws.Range('D13').Value = '3.538.87'
$ws.Range('E13').Value = '  +0.84%  '
$ws.Range('D14').Value = "'18.22"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.63%  '
$ws.Range('D15').Value = "'7.63"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.36%  '
$ws.Range('D16').Value = '3.017.16'
$ws.Range('E16').Value = '  -0.63%  '
$ws.Range('E17').Value = '  +1.99%  '
$ws.Range('D18').Value = "'10.53"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.47%  '
$ws.Range('D19').Value = '51.100.43'
$ws.Range('E19').Value = '  -1.09%  '
$ws.Range('D20').Value = "'3.16"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.78%  '
$ws.Range('E21').Value = '  -2.21%  '
$ws.Range('E22').Value = '  -1.16%  '
$ws.Range('D23').Value = "'69.54"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.63%  '
$ws.Range('D24').Value = "'263.25"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.56%  '
$ws.Range('D25').Value = "'3.13"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.37%  '
$ws.Range('D26').Value = "'7.85"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -7.05%  '
$ws.Range('D27').Value = "'26.68"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.65%  '
$ws.Range('E28').Value = '  -0.17%  '
$ws.Range('D29').Value = "'7.09"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.29%  '
$ws.Range('E30').Value = '  -5.99%  '
$ws.Range('E31').Value = '  -3.10%  '
$ws.Range('D32').Value = "'10.47"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.89%  '
$ws.Range('D33').Value = "'0.0488"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +9.06%  '
$ws.Range('D34').Value = "'35.51"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.03%  '
$ws.Range('D35').Value = "'2.06"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.65%  '
$ws.Range('D36').Value = "'49.96"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.23%  '
$ws.Range('E37').Value = '  +0.00%  '
$ws.Range('E38').Value = '  +0.69%  '
$ws.Range('E39').Value = '  -1.58%  '
$ws.Range('D40').Value = "'129.14"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.17%  '
$ws.Range('E41').Value = '  -3.41%  '
$ws.Range('E42').Value = '  -2.71%  '
$ws.Range('E43').Value = '  -1.80%  '
$ws.Range('E44').Value = '  +1.52%  '
$ws.Range('E45').Value = '  -2.34%  '
$ws.Range('D46').Value = "'21.61"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.18%  '
$ws.Range('E47').Value = '  +2.78%  '
$ws.Range('E48').Value = '  -0.32%  '
$ws.Range('D49').Value = '2.061.92'
$ws.Range('E49').Value = '  +1.64%  '
$ws.Range('D50').Value = "'0.0319"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.39%  '
$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').Value = "'5.39"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.19%  '
